# Europa.xlsx edit: shift Sheet1's data block two columns to the right
# (A:D -> C:F) and update the active-sheet/selection UI state to match.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Sheets.Item("Sheet1")
$ws2 = $wb.Sheets.Item("Sheet2")
$ws3 = $wb.Sheets.Item("Sheet3")

# --- Sheet1: insert two blank columns before column A, pushing the whole
#     A2:D12 block (values, shared-string refs, number-format styles, and
#     the custom column width) to C2:F12 in one structural move. ---
$ws1.Columns("A:B").Insert()

# Select column A (the now-empty leading column) as the sheet's selection,
# with the view anchored at A1.
$ws1.Range("A1:A1048576").Select()

# --- Sheet2: selection moves from F4 to F7 ---
$ws2.Range("F7").Select()

# --- Sheet3: selection moves from I11 to F8 ---
$ws3.Range("F8").Select()

# --- Sheet1 becomes the active/visible tab (previously Sheet3) ---
$ws1.Activate()
